# Refresh the cryptocurrency price/volume snapshot (and two rank-tie swaps)
# to match the latest scrape, per the GitHub Actions commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text happens to look like a plain number (e.g. "418.69")
# are written with a leading apostrophe - Excel's normal 'force text' marker -
# so they stay text cells (matching their neighbours) instead of being
# auto-converted to numeric values.

$ws.Range("D2").Value = "66.595.49"
$ws.Range("E2").Value = "  +3.82%  "

$ws.Range("D3").Value = "3.691.82"
$ws.Range("E3").Value = "  +5.96%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").Value = "`'418.69"
$ws.Range("E5").Value = "  +0.24%  "

$ws.Range("D6").Value = "`'130.39"
$ws.Range("E6").Value = "  +0.67%  "

$ws.Range("D7").Value = "3.684.03"
$ws.Range("E7").Value = "  +6.08%  "

$ws.Range("D8").Value = "`'0.643"
$ws.Range("E8").Value = "  +0.40%  "

$ws.Range("E9").Value = "  -0.06%  "

$ws.Range("D10").Value = "`'0.760"
$ws.Range("E10").Value = "  -3.44%  "

$ws.Range("E11").Value = "  +11.18%  "

$ws.Range("D12").Value = "`'0.0000394"
$ws.Range("E12").Value = "  +47.85%  "

$ws.Range("D13").Value = "`'42.78"
$ws.Range("E13").Value = "  -0.60%  "

$ws.Range("D14").Value = "`'10.60"
$ws.Range("E14").Value = "  +7.70%  "

$ws.Range("D15").Value = "4.278.06"
$ws.Range("E15").Value = "  +6.12%  "

$ws.Range("E16").Value = "  -0.48%  "

$ws.Range("B17").Value = "Chainlink"
$ws.Range("C17").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D17").Value = "`'20.54"
$ws.Range("E17").Value = "  +0.55%  "

$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.716.29"
$ws.Range("E18").Value = "  +7.06%  "

$ws.Range("E19").Value = "  +6.11%  "

$ws.Range("E20").Value = "  +2.38%  "

$ws.Range("D21").Value = "66.635.36"
$ws.Range("E21").Value = "  +4.12%  "

$ws.Range("D22").Value = "`'443.42"
$ws.Range("E22").Value = "  -3.32%  "

$ws.Range("D23").Value = "`'16.28"
$ws.Range("E23").Value = "  +20.53%  "

$ws.Range("D24").Value = "`'89.90"
$ws.Range("E24").Value = "  -1.05%  "

$ws.Range("D25").Value = "`'3.13"
$ws.Range("E25").Value = "  -4.38%  "

$ws.Range("E26").Value = "  +9.87%  "

$ws.Range("D27").Value = "`'10.25"
$ws.Range("E27").Value = "  +0.06%  "

$ws.Range("E28").Value = "  -0.75%  "

$ws.Range("E29").Value = "  +4.23%  "

$ws.Range("D30").Value = "`'12.65"
$ws.Range("E30").Value = "  +0.59%  "

$ws.Range("E31").Value = "  +6.68%  "

$ws.Range("D32").Value = "`'2.70"
$ws.Range("E32").Value = "  -0.05%  "

$ws.Range("D33").Value = "`'7.26"
$ws.Range("E33").Value = "  -4.76%  "

$ws.Range("E34").Value = "  -0.89%  "

$ws.Range("D35").Value = "`'41.17"
$ws.Range("E35").Value = "  +2.40%  "

$ws.Range("D36").Value = "`'57.33"
$ws.Range("E36").Value = "  -1.26%  "

$ws.Range("D37").Value = "`'1.00"
$ws.Range("E37").Value = "  -0.10%  "

$ws.Range("D38").Value = "`'0.0492"
$ws.Range("E38").Value = "  -4.72%  "

$ws.Range("D39").Value = "`'3.22"
$ws.Range("E39").Value = "  +37.97%  "

$ws.Range("D40").Value = "0.0₃0743"
$ws.Range("E40").Value = "  +12.71%  "

$ws.Range("E41").Value = "  +3.42%  "

$ws.Range("D42").Value = "`'28.51"
$ws.Range("E42").Value = "  +29.94%  "

$ws.Range("D43").Value = "`'0.997"
$ws.Range("E43").Value = "  -0.10%  "

$ws.Range("E44").Value = "  +1.19%  "

$ws.Range("D45").Value = "`'148.15"
$ws.Range("E45").Value = "  +2.51%  "

$ws.Range("E46").Value = "  +3.48%  "

$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D47").Value = "`'4.37"
$ws.Range("E47").Value = "  -2.95%  "

$ws.Range("B48").Value = "Stacks"
$ws.Range("C48").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D48").Value = "`'2.89"
$ws.Range("E48").Value = "  -7.40%  "

$ws.Range("E49").Value = "  -4.74%  "

$ws.Range("D50").Value = "`'2.56"
$ws.Range("E50").Value = "  -6.60%  "

$ws.Range("E51").Value = "  +15.57%  "
